$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 84 (pushes existing rows 84-131 down to 85-132,
# carrying over the row-84 "D" column date style to the new row).
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new weekly record.
$ws.Cells.Item(84, 1).Value = 7
$ws.Cells.Item(84, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(84, 3).Value = "Ñuble"
$ws.Cells.Item(84, 4).Value = 44960
$ws.Cells.Item(84, 5).Value = 16
$ws.Cells.Item(84, 6).Value = 100112021
$ws.Cells.Item(84, 7).Value = "Ají"
$ws.Cells.Item(84, 8).Value = "Americana (o)"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 60
$ws.Cells.Item(84, 11).Value = 11000
$ws.Cells.Item(84, 12).Value = 12000
$ws.Cells.Item(84, 13).Value = 11500
$ws.Cells.Item(84, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(84, 15).Value = "Región del Maule"
$ws.Cells.Item(84, 16).Value = 767
$ws.Cells.Item(84, 17).Value = 15
$ws.Cells.Item(84, 18).Value = "Hortaliza"
